$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("展览")
$ws.Cells.Item(4, 6).Value = 656
$ws.Cells.Item(5, 6).Value = 2952
$ws.Cells.Item(10, 6).Value = 6931
$ws.Cells.Item(11, 6).Value = 41
$ws.Cells.Item(12, 6).Value = 85
$ws.Cells.Item(14, 6).Value = 605
$ws.Cells.Item(15, 6).Value = 1496
$ws.Cells.Item(16, 6).Value = 1120
$ws.Cells.Item(17, 6).Value = 2244
$ws.Cells.Item(18, 6).Value = 1492
$ws.Cells.Item(20, 6).Value = 121
$ws.Cells.Item(21, 6).Value = 1114
$ws.Cells.Item(22, 6).Value = 131
$ws.Cells.Item(23, 6).Value = 180
$ws.Cells.Item(24, 6).Value = 347
$ws.Cells.Item(25, 6).Value = 13
$ws.Cells.Item(26, 6).Value = 1732
$ws.Cells.Item(27, 6).Value = 1697
$ws.Cells.Item(28, 6).Value = 1033
$ws.Cells.Item(30, 6).Value = 1667
$ws.Cells.Item(31, 6).Value = 1226
$ws.Cells.Item(32, 6).Value = 141
$ws.Cells.Item(33, 6).Value = 588
$ws.Cells.Item(35, 6).Value = 430
$ws.Cells.Item(36, 6).Value = 17
$ws.Cells.Item(37, 6).Value = 2483
$ws.Cells.Item(38, 6).Value = 2732
$ws.Cells.Item(39, 6).Value = 74
$ws.Cells.Item(40, 6).Value = 10
$ws.Cells.Item(41, 6).Value = 186
$ws.Cells.Item(42, 6).Value = 17
$ws.Cells.Item(43, 6).Value = 29
$ws.Cells.Item(44, 6).Value = 318
$ws.Cells.Item(45, 6).Value = 126
$ws.Cells.Item(46, 6).Value = 170
$ws.Cells.Item(47, 6).Value = 158
$ws.Cells.Item(48, 6).Value = 10

$ws = $wb.Worksheets.Item("演出")
$ws.Cells.Item(7, 6).Value = 168
$ws.Cells.Item(10, 6).Value = 32
$ws.Cells.Item(22, 6).Value = 340
$ws.Cells.Item(23, 6).Value = 476
$ws.Cells.Item(29, 6).Value = 22
$ws.Cells.Item(30, 6).Value = 26

$ws = $wb.Worksheets.Item("本地生活")
$ws.Cells.Item(4, 6).Value = 542
$ws.Cells.Item(6, 6).Value = 1687
$ws.Cells.Item(7, 6).Value = 1854
$ws.Cells.Item(8, 6).Value = 2733
$ws.Cells.Item(9, 6).Value = 1020
$ws.Cells.Item(10, 6).Value = 934
$ws.Cells.Item(12, 6).Value = 272
$ws.Cells.Item(13, 6).Value = 1479
$ws.Cells.Item(14, 6).Value = 7364

$ws = $wb.Worksheets.Item("全部类型")
$ws.Cells.Item(3, 6).Value = 542
$ws.Cells.Item(4, 6).Value = 656
$ws.Cells.Item(5, 6).Value = 2952
$ws.Cells.Item(6, 6).Value = 1687
$ws.Cells.Item(8, 6).Value = 2733
$ws.Cells.Item(9, 6).Value = 6931
$ws.Cells.Item(10, 6).Value = 1020
$ws.Cells.Item(11, 6).Value = 41
$ws.Cells.Item(13, 6).Value = 168
$ws.Cells.Item(14, 6).Value = 272
$ws.Cells.Item(15, 6).Value = 1120
$ws.Cells.Item(16, 6).Value = 2244
$ws.Cells.Item(17, 6).Value = 1492
$ws.Cells.Item(18, 6).Value = 122
$ws.Cells.Item(20, 6).Value = 1114
$ws.Cells.Item(22, 6).Value = 1732
$ws.Cells.Item(25, 6).Value = 1667
$ws.Cells.Item(26, 6).Value = 1226
$ws.Cells.Item(27, 6).Value = 141
$ws.Cells.Item(29, 6).Value = 588
$ws.Cells.Item(33, 6).Value = 340
$ws.Cells.Item(34, 6).Value = 476
$ws.Cells.Item(35, 6).Value = 430
$ws.Cells.Item(37, 6).Value = 17
$ws.Cells.Item(38, 6).Value = 2484
$ws.Cells.Item(39, 6).Value = 2732
$ws.Cells.Item(40, 6).Value = 74
$ws.Cells.Item(41, 6).Value = 186
$ws.Cells.Item(42, 6).Value = 17
$ws.Cells.Item(43, 6).Value = 29
$ws.Cells.Item(44, 6).Value = 318
$ws.Cells.Item(45, 6).Value = 126
$ws.Cells.Item(46, 6).Value = 170

